$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 181
$ws.Range("I8").Value = 78.28570999999999
$ws.Range("K8").Value = 234.85713
$ws.Range("M8").Value = -95.85712999999998

$ws.Range("H17").Value = 1353.4375
$ws.Range("J17").Value = 1353.4375
$ws.Range("L17").Value = 4060.3125
$ws.Range("N17").Value = -4396.3125

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H46").Value = 5209321
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5209321
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 15627963
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -15628201

$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 9942.929
$ws.Range("J51").Value = 3094.6843
$ws.Range("K51").Value = 9942.929
$ws.Range("L51").Value = 3094.6843
$ws.Range("M51").Value = -9458.929
$ws.Range("N51").Value = -4062.6843

$ws.Range("H60").Value = 5209321
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 5209321
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 15627963
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -15628931

$ws.Range("H103").Value = 1440.1765
$ws.Range("J103").Value = 1234.9
$ws.Range("L103").Value = 3704.7
$ws.Range("N103").Value = -4876.700000000001

$ws.Range("H138").Value = 3040.611
$ws.Range("I138").Value = 2181.5789
$ws.Range("J138").Value = 3270.493
$ws.Range("K138").Value = 6544.736699999999
$ws.Range("L138").Value = 9811.478999999999
$ws.Range("M138").Value = -1404.736699999999
$ws.Range("N138").Value = -20091.479

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 408.16666
$ws.Range("I5").Value = 636.6667
$ws.Range("J5").Value = 179.66667
$ws.Range("K5").Value = 636.6667
$ws.Range("L5").Value = 179.66667
$ws.Range("M5").Value = -524.6667
$ws.Range("N5").Value = -403.66667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 408.16666
$ws.Range("I4").Value = 636.6667
$ws.Range("J4").Value = 179.66667
$ws.Range("K4").Value = 636.6667
$ws.Range("L4").Value = 179.66667
$ws.Range("M4").Value = -521.6667
$ws.Range("N4").Value = -409.66667

$ws.Range("H22").Value = 186.5
$ws.Range("I22").Value = 191.71428
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 191.71428
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = -18.71428
$ws.Range("N22").Value = -496

$ws.Range("H105").Value = 401586.2
$ws.Range("I105").Value = 501990
$ws.Range("J105").Value = 334650.34
$ws.Range("K105").Value = 501990
$ws.Range("L105").Value = 334650.34
$ws.Range("M105").Value = -500243
$ws.Range("N105").Value = -338144.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 532.7
$ws.Range("I22").Value = 145.2
$ws.Range("J22").Value = 920.2
$ws.Range("K22").Value = 145.2
$ws.Range("L22").Value = 920.2
$ws.Range("M22").Value = 204.8
$ws.Range("N22").Value = -1620.2

$ws.Range("H70").Value = 11999.25
$ws.Range("J70").Value = 11999.25
$ws.Range("L70").Value = 11999.25
$ws.Range("N70").Value = -12629.25

$ws.Range("H73").Value = 11999.25
$ws.Range("J73").Value = 11999.25
$ws.Range("L73").Value = 11999.25
$ws.Range("N73").Value = -14183.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 706.6818
$ws.Range("I122").Value = 539.7143
$ws.Range("J122").Value = 784.6
$ws.Range("K122").Value = 4857.428699999999
$ws.Range("L122").Value = 7061.400000000001
$ws.Range("M122").Value = -2407.428699999999
$ws.Range("N122").Value = -11961.4

$ws.Range("H131").Value = 834.15
$ws.Range("I131").Value = 570.125
$ws.Range("J131").Value = 884.4405
$ws.Range("K131").Value = 1710.375
$ws.Range("L131").Value = 2653.3215
$ws.Range("M131").Value = 3329.625
$ws.Range("N131").Value = -12733.3215

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 40249.5
$ws.Range("J51").Value = 40249.5
$ws.Range("L51").Value = 40249.5
$ws.Range("N51").Value = -41267.5

$ws.Range("H62").Value = 15000
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372

$ws.Range("H65").Value = 15000
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864

$ws.Range("H74").Value = 8054
$ws.Range("J74").Value = 8054
$ws.Range("L74").Value = 8054
$ws.Range("N74").Value = -9926

$ws.Range("H77").Value = 8054
$ws.Range("J77").Value = 8054
$ws.Range("L77").Value = 24162
$ws.Range("N77").Value = -33522

$ws.Range("H132").Value = 2739.8823
$ws.Range("I132").Value = 1950.7142
$ws.Range("K132").Value = 5852.142599999999
$ws.Range("M132").Value = -3322.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1164.7059
$ws.Range("I22").Value = 3466.6667
$ws.Range("J22").Value = 671.4286
$ws.Range("K22").Value = 3466.6667
$ws.Range("L22").Value = 671.4286
$ws.Range("M22").Value = -3171.6667
$ws.Range("N22").Value = -1261.4286

$ws.Range("H27").Value = 1164.7059
$ws.Range("I27").Value = 3466.6667
$ws.Range("J27").Value = 671.4286
$ws.Range("K27").Value = 3466.6667
$ws.Range("L27").Value = 671.4286
$ws.Range("M27").Value = -3359.6667
$ws.Range("N27").Value = -885.4286

$ws.Range("H46").Value = 6227.9
$ws.Range("I46").Value = 4395
$ws.Range("J46").Value = 7449.8335
$ws.Range("K46").Value = 4395
$ws.Range("L46").Value = 7449.8335
$ws.Range("M46").Value = -4207
$ws.Range("N46").Value = -7825.8335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 13677.777
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 13677.777
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 13677.777
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -14263.777

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
